# 107_2.xlsx confirmations sheet: the per-branch breakdown rows keep their
# original generic labels ("New nominations", "Carryover nominations", ...)
# but the diff renames them to include the branch name (e.g. "Civilian, New
# nominations"), drops the standalone "Summary" section header row, and
# renames / reorders the two grand totals that used to sit right under it
# ("carried over" <-> "received"/"new").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the "Summary" section header (row 36); everything below it shifts
#    up by one row (old row 37 -> new row 36, ... old row 42 -> new row 41).
$ws.Rows.Item(36).Delete() | Out-Null

# 2) The old "Total nominations carried over from the First Session " (791)
#    and "Total nominations Received this Session " (23045) totals swap
#    places & identity: the new row 36 ("Total new nominations") must hold
#    23045 (with its thousands-separator number format) and the new row 37
#    ("Total carryover nominations") must hold 791 (plain number format).
#    Swap values+formats via a scratch cell well outside the sheet's used
#    range so nothing stray is left behind once it's cleared.
$scratch = $ws.Cells.Item(1000, 26)
$ws.Cells.Item(36, 2).Copy($scratch) | Out-Null
$ws.Cells.Item(37, 2).Copy($ws.Cells.Item(36, 2)) | Out-Null
$scratch.Copy($ws.Cells.Item(37, 2)) | Out-Null
$scratch.Clear() | Out-Null

# 3) Relabel column A, top to bottom, to match the new wording.
$ws.Cells.Item(7, 1).Value = '     Civilian, New nominations'
$ws.Cells.Item(8, 1).Value = '     Civilian, Carryover nominations'
$ws.Cells.Item(9, 1).Value = '     Civilian, Confirmed '
$ws.Cells.Item(10, 1).Value = '     Civilian, Withdrawn '
$ws.Cells.Item(11, 1).Value = '     Civilian, Returned to White House '

$ws.Cells.Item(13, 1).Value = '     Other Civilian, New nominations'
$ws.Cells.Item(14, 1).Value = '     Other Civilian, Carryover nominations'
$ws.Cells.Item(15, 1).Value = '     Other Civilian, Confirmed '
$ws.Cells.Item(16, 1).Value = '     Other Civilian, Returned to White House '

$ws.Cells.Item(18, 1).Value = '     Air Force, New nominations'
$ws.Cells.Item(19, 1).Value = '     Air Force, Carryover nominations'
$ws.Cells.Item(20, 1).Value = '     Air Force, Confirmed '
$ws.Cells.Item(21, 1).Value = '     Air Force, Returned to White House '

$ws.Cells.Item(23, 1).Value = '     Army, New nominations'
$ws.Cells.Item(24, 1).Value = '     Army, Carryover nominations'
$ws.Cells.Item(25, 1).Value = '     Army, Confirmed '
$ws.Cells.Item(26, 1).Value = '     Army, Returned to White House '

$ws.Cells.Item(28, 1).Value = '     Navy, New nominations'
$ws.Cells.Item(29, 1).Value = '     Navy, Confirmed '
$ws.Cells.Item(30, 1).Value = '     Navy, Returned to White House '

$ws.Cells.Item(32, 1).Value = '     Marine Corps, New nominations'
$ws.Cells.Item(33, 1).Value = '     Marine Corps, Carryover nominations'
$ws.Cells.Item(34, 1).Value = '     Marine Corps, Confirmed '
$ws.Cells.Item(35, 1).Value = '     Marine Corps, Returned to White House '

$ws.Cells.Item(36, 1).Value = 'Total new nominations'
$ws.Cells.Item(37, 1).Value = 'Total carryover nominations'
$ws.Cells.Item(38, 1).Value = 'Total confirmed '
$ws.Cells.Item(39, 1).Value = 'Total unconfirmed '
$ws.Cells.Item(40, 1).Value = 'Total withdrawn '
$ws.Cells.Item(41, 1).Value = 'Total returned'
